$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/7f1fe05bd94b45868680f631bed471ca2a2e016c"

$zhFile = "4430e4f1-aebe-483e-be79-9220610db35f.cdac7166739e2257d9b8ca96ddaaf551f81d709f.zh-cn.xlf"
$deFile = "4430e4f1-aebe-483e-be79-9220610db35f.cdac7166739e2257d9b8ca96ddaaf551f81d709f.de-de.xlf"

# ---------------------------------------------------------------------------
# 1. "Handoff transform failed" -> "Ready for handoff" everywhere it appears
#    (Overview!B2, Overview!C2, zh-cn!B2, de-de!B2 all share this string)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = "Ready for handoff"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: add the Latest Handoff File link/date, flip Handoff Reason
# ---------------------------------------------------------------------------
$wsZh.Range("D2").Value = "2016-01-14 05:38:24"
$wsZh.Range("H2").Value = "Include"

# Add the new Latest Handoff File link on C2. A2 and A3 already carry the
# correct hyperlinks from the original workbook, so leave them alone and
# only wire up the new cell.
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "$repoBase/e2e/$zhFile", "", "", $zhFile)

# Match the look of the other hyperlink cells (underlined, hyperlink blue)
$wsZh.Range("C2").Font.Underline = 2
$wsZh.Range("C2").Font.Color = 15570276

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of edit, different file name / timestamp
# ---------------------------------------------------------------------------
$wsDe.Range("D2").Value = "2016-01-14 05:38:47"
$wsDe.Range("H2").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/e2e/4430e4f1-aebe-483e-be79-9220610db35f.md", "", "", "4430e4f1-aebe-483e-be79-9220610db35f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$repoBase/e2e/$deFile", "", "", $deFile)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/.localization-config", "", "", ".localization-config")

foreach ($addr in @("A2", "C2", "A3")) {
    $wsDe.Range($addr).Font.Underline = 2
    $wsDe.Range($addr).Font.Color = 15570276
}
